$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with the default (unstyled) data format, used to restore
# D-column cells to their original style after a text-coercion trick below.
$defaultStyle = $ws.Range("D22").Style

# Row 2
$ws.Range("D2").Value = "49.534.45"
$ws.Range("E2").Value = "  -0.74%  "

# Row 3
$ws.Range("D3").Value = "2.632.84"
$ws.Range("E3").Value = "  -0.88%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "111.40"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  -1.76%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "325.48"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  -0.79%  "

# Row 7
$ws.Range("E7").Value = "  -1.39%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("E9").Value = "  -1.83%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.44"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = "  -4.40%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.14"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = "  -0.54%  "

# Row 12
$ws.Range("E12").Value = "  -1.88%  "

# Row 13
$ws.Range("E13").Value = "  +1.26%  "

# Row 14
$ws.Range("E14").Value = "  +0.01%  "

# Row 15
$ws.Range("D15").Value = "3.039.13"
$ws.Range("E15").Value = "  +1.09%  "

# Row 16
$ws.Range("D16").Value = "2.628.45"
$ws.Range("E16").Value = "  -0.73%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.851"
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = "  -2.87%  "

# Row 18
$ws.Range("D18").Value = "49.445.97"
$ws.Range("E18").Value = "  -0.76%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.04"
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = "  -0.80%  "

# Row 20
$ws.Range("E20").Value = "  -1.52%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.90"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = "  -1.43%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "267.83"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = "  -3.30%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.93"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "  -4.72%  "

# Row 25
$ws.Range("E25").Value = "  -1.99%  "

# Row 26
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.01"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = "  -3.27%  "

# Row 27
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = "  +0.13%  "

# Row 28
$ws.Range("E28").Value = "  +1.95%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.19"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = "  -1.09%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.139"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = "  -1.67%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.53"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = "  -4.70%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.58"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = "  -1.39%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.47"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = "  +0.84%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0809"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = "  +0.98%  "

# Row 35
$ws.Range("E35").Value = "  -0.07%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.99"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = "  -2.87%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.94"
$ws.Range("D37").Style = $defaultStyle
$ws.Range("E37").Value = "  +3.62%  "

# Row 38
$ws.Range("E38").Value = "  -2.79%  "

# Row 39
$ws.Range("E39").Value = "  -0.08%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "128.83"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  +2.67%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.03"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "  +3.30%  "

# Row 42
$ws.Range("E42").Value = "  -1.62%  "

# Row 43
$ws.Range("E43").Value = "  -1.01%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0327"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = "  +3.53%  "

# Row 45
$ws.Range("D45").Value = "2.042.10"
$ws.Range("E45").Value = "  -1.43%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.17"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = "  +8.95%  "

# Row 47
$ws.Range("E47").Value = "  -3.75%  "

# Row 48
$ws.Range("E48").Value = "  -4.38%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.84"
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Value = "  -3.56%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.21"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = "  -3.65%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "58.49"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = "  +1.24%  "
